$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Formula = "=B2+C2+D2+E2+F2+G2+H2+I2+J2+K2"
$ws.Range("L3:L22").Formula = "=B3+C3+D3+E3+F3+G3+H3+I3+J3+K3"

$ws.Range("P18").Select()
